$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.442.97"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +3.23%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'2.307.16"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +2.17%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.19%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'310.15"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +0.89%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'104.62"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +7.90%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.531"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +1.17%  "
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "'  +0.13%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "'  +8.10%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'36.33"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +4.73%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'52.49"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +0.75%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.0813"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  -0.20%  "
$ws.Range("E12").ClearFormats()
$ws.Range("E13").Value = "'  -0.57%  "
$ws.Range("E13").ClearFormats()
$ws.Range("E14").Value = "'  +2.47%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'2.671.96"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +2.44%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'15.09"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +3.62%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'2.318.69"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +2.70%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'0.805"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +2.77%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'43.384.28"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +3.40%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'12.00"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -1.78%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'0.0₃0924"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +2.63%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'6.16"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +4.02%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'68.06"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +0.93%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'240.89"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +2.34%  "
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = "'  +3.15%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'2.61"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +1.78%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'0.997"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  -0.23%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'24.90"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +6.17%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'2.24"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +5.43%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'36.58"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  -0.62%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'9.59"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  +0.72%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'164.17"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  -0.08%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'5.24"
$ws.Range("D33").ClearFormats()
$ws.Range("E34").Value = "'  +0.09%  "
$ws.Range("E34").ClearFormats()
$ws.Range("E35").Value = "'  +4.73%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'2.54"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +6.91%  "
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = "'  +1.61%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'3.02"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -1.80%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'4.51"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +9.36%  "
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = "'  +4.41%  "
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'  +2.70%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.115"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +0.57%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'2.63"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +17.84%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'1.988.03"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +2.46%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.0290"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +3.24%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'18.95"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +2.53%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'3.07"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +5.05%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'10.12"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +4.43%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'58.11"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +7.74%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'1.59"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +8.65%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'2.91"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +1.32%  "
$ws.Range("E51").ClearFormats()
